$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "-"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "-"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("B16").Value = "-"
